# chore: update Sheets via scheduled runner
# Refreshes cached market-board figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leves across the ALC/ARM/BSM/CUL/GSM/LTW/WVR
# sheets, matching newer Universalis pricing data.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 120
$ws.Range("H120").Value = 28000
$ws.Range("J120").Value = 28000
$ws.Range("L120").Value = 28000
$ws.Range("N120").Value = -37676

# Row 129
$ws.Range("H129").Value = 1138.3055
$ws.Range("I129").Value = 319.83334
$ws.Range("J129").Value = 1302
$ws.Range("K129").Value = 959.5000200000001
$ws.Range("L129").Value = 3906
$ws.Range("M129").Value = 4040.49998
$ws.Range("N129").Value = -13906


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 45
$ws.Range("H45").Value = 2084.5789
$ws.Range("I45").Value = 2252.7693
$ws.Range("J45").Value = 1720.1666
$ws.Range("K45").Value = 2252.7693
$ws.Range("L45").Value = 1720.1666
$ws.Range("M45").Value = -1875.7693
$ws.Range("N45").Value = -2474.1666

# Row 121
$ws.Range("H121").Value = 35000
$ws.Range("J121").Value = 35000
$ws.Range("L121").Value = 35000
$ws.Range("N121").Value = -38494

# Row 122
$ws.Range("H122").Value = 2101.2
$ws.Range("I122").Value = 1436
$ws.Range("J122").Value = 3099
$ws.Range("K122").Value = 4308
$ws.Range("L122").Value = 9297
$ws.Range("M122").Value = -1858
$ws.Range("N122").Value = -14197

# Row 123
$ws.Range("H123").Value = 60732.25
$ws.Range("J123").Value = 60732.25
$ws.Range("L123").Value = 60732.25
$ws.Range("N123").Value = -70532.25


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 26
$ws.Range("H26").Value = 23203.8
$ws.Range("I26").Value = 20833.334
$ws.Range("J26").Value = 26759.5
$ws.Range("K26").Value = 20833.334
$ws.Range("L26").Value = 26759.5
$ws.Range("M26").Value = -20541.334
$ws.Range("N26").Value = -27343.5

# Row 96
$ws.Range("H96").Value = 20424.666
$ws.Range("I96").Value = 9289.25
$ws.Range("J96").Value = 25992.375
$ws.Range("K96").Value = 9289.25
$ws.Range("L96").Value = 25992.375
$ws.Range("M96").Value = -6543.25
$ws.Range("N96").Value = -31484.375

# Row 99
$ws.Range("H99").Value = 1080.826
$ws.Range("I99").Value = 1097.5883
$ws.Range("K99").Value = 1097.5883
$ws.Range("M99").Value = 400.4117000000001


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 31
$ws.Range("H31").Value = 2483.3333
$ws.Range("J31").Value = 2483.3333
$ws.Range("L31").Value = 7449.999899999999
$ws.Range("N31").Value = -8025.999899999999

# Row 49
$ws.Range("H49").Value = 1558.3334
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1558.3334
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 4675.0002
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -4987.0002

# Row 63
$ws.Range("H63").Value = 11577.889
$ws.Range("I63").Value = 8475
$ws.Range("J63").Value = 14060.2
$ws.Range("K63").Value = 25425
$ws.Range("L63").Value = 42180.60000000001
$ws.Range("M63").Value = -24676
$ws.Range("N63").Value = -43678.60000000001

# Row 64
$ws.Range("H64").Value = 4443.5654
$ws.Range("I64").Value = 1537.3334
$ws.Range("J64").Value = 4879.5
$ws.Range("K64").Value = 4612.0002
$ws.Range("L64").Value = 14638.5
$ws.Range("M64").Value = -4342.0002
$ws.Range("N64").Value = -15178.5

# Row 66
$ws.Range("H66").Value = 11577.889
$ws.Range("I66").Value = 8475
$ws.Range("J66").Value = 14060.2
$ws.Range("K66").Value = 76275
$ws.Range("L66").Value = 126541.8
$ws.Range("M66").Value = -72531
$ws.Range("N66").Value = -134029.8

# Row 67
$ws.Range("H67").Value = 4443.5654
$ws.Range("I67").Value = 1537.3334
$ws.Range("J67").Value = 4879.5
$ws.Range("K67").Value = 4612.0002
$ws.Range("L67").Value = 14638.5
$ws.Range("M67").Value = -3676.0002
$ws.Range("N67").Value = -16510.5

# Row 74
$ws.Range("H74").Value = 14852.083
$ws.Range("I74").Value = 4113
$ws.Range("J74").Value = 16999.9
$ws.Range("K74").Value = 12339
$ws.Range("L74").Value = 50999.7
$ws.Range("M74").Value = -11278
$ws.Range("N74").Value = -53121.7

# Row 77
$ws.Range("H77").Value = 14852.083
$ws.Range("I77").Value = 4113
$ws.Range("J77").Value = 16999.9
$ws.Range("K77").Value = 37017
$ws.Range("L77").Value = 152999.1
$ws.Range("M77").Value = -31713
$ws.Range("N77").Value = -163607.1

# Row 92
$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -3996

# Row 94
$ws.Range("H94").Value = 4691.3335
$ws.Range("J94").Value = 4659.625
$ws.Range("L94").Value = 13978.875
$ws.Range("N94").Value = -15330.875

# Row 97
$ws.Range("H97").Value = 1127.5
$ws.Range("I97").Value = 836.6667
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 2510.0001
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -2014.0001
$ws.Range("N97").Value = -6992

# Row 98
$ws.Range("H98").Value = 436
$ws.Range("I98").Value = 436
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1308
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 190
$ws.Range("N98").ClearContents()

# Row 99
$ws.Range("H99").Value = 6939.091
$ws.Range("I99").Value = 3500
$ws.Range("K99").Value = 10500
$ws.Range("M99").Value = -8254

# Row 100
$ws.Range("H100").Value = 15114.214
$ws.Range("I100").Value = 13999
$ws.Range("J100").Value = 15200
$ws.Range("K100").Value = 41997
$ws.Range("L100").Value = 45600
$ws.Range("M100").Value = -41186
$ws.Range("N100").Value = -47222

# Row 101
$ws.Range("H101").Value = 15725
$ws.Range("J101").Value = 15725
$ws.Range("L101").Value = 47175
$ws.Range("N101").Value = -52043

# Row 114
$ws.Range("H114").Value = 834.4211
$ws.Range("I114").Value = 941.1539
$ws.Range("J114").Value = 603.1667
$ws.Range("K114").Value = 2823.4617
$ws.Range("L114").Value = 1809.5001
$ws.Range("M114").Value = 430.5383000000002
$ws.Range("N114").Value = -8317.500099999999

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# Row 130
$ws.Range("H130").Value = 1360
$ws.Range("I130").Value = 1000
$ws.Range("K130").Value = 3000
$ws.Range("M130").Value = 2020

# Row 131
$ws.Range("H131").Value = 764.5833
$ws.Range("J131").Value = 981.4286
$ws.Range("L131").Value = 2944.2858
$ws.Range("N131").Value = -13024.2858


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 97
$ws.Range("H97").Value = 1817.5264
$ws.Range("I97").Value = 1710.909
$ws.Range("J97").Value = 1964.125
$ws.Range("K97").Value = 1710.909
$ws.Range("L97").Value = 1964.125
$ws.Range("M97").Value = -1214.909
$ws.Range("N97").Value = -2956.125

# Row 122
$ws.Range("H122").Value = 3140.8538
$ws.Range("I122").Value = 2735.5417
$ws.Range("J122").Value = 3713.0588
$ws.Range("K122").Value = 8206.625100000001
$ws.Range("L122").Value = 11139.1764
$ws.Range("M122").Value = -5756.625100000001
$ws.Range("N122").Value = -16039.1764


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 697.1667
$ws.Range("I22").Value = 483.26666
$ws.Range("J22").Value = 1766.6666
$ws.Range("K22").Value = 483.26666
$ws.Range("L22").Value = 1766.6666
$ws.Range("M22").Value = -188.26666
$ws.Range("N22").Value = -2356.6666

# Row 27
$ws.Range("H27").Value = 697.1667
$ws.Range("I27").Value = 483.26666
$ws.Range("J27").Value = 1766.6666
$ws.Range("K27").Value = 483.26666
$ws.Range("L27").Value = 1766.6666
$ws.Range("M27").Value = -376.26666
$ws.Range("N27").Value = -1980.6666

# Row 122
$ws.Range("H122").Value = 3182.5454
$ws.Range("I122").Value = 3560.375
$ws.Range("J122").Value = 2175
$ws.Range("K122").Value = 10681.125
$ws.Range("L122").Value = 6525
$ws.Range("M122").Value = -8231.125
$ws.Range("N122").Value = -11425


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 96
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 373
$ws.Range("N96").ClearContents()

